# Avoid recalculating the SUM() column (B) when D6/E6 switch from numeric
# durations to their text "h:mm" representations -- the source edit was a
# content-only change and the cached formula results were left untouched.
$excel.Calculation = -4135  # xlCalculationManual

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "total tea units sold" column header to the new Chai-specific
# wording. This also renames the matching column in the worksheet Table.
$ws.Range("B1").Value = "印度奶茶总销售数量（件）"

# Row 6 (2024-05-04 ish) previously stored raw sales-count numbers in D/E;
# replace them with elapsed-time style text values.
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "17:05"
